$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 399, shifting rows 399:475 down to 400:476
$ws.Rows("399:399").Insert()

# Populate the new row 399 with the new price record
$ws.Cells.Item(399, 1).Value = 4
$ws.Cells.Item(399, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(399, 3).Value = "Los Lagos"
$ws.Cells.Item(399, 4).Value = 45211
$ws.Cells.Item(399, 5).Value = 10
$ws.Cells.Item(399, 6).Value = "Fruta"
$ws.Cells.Item(399, 7).Value = 100108
$ws.Cells.Item(399, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(399, 9).Value = 100108005
$ws.Cells.Item(399, 10).Value = "Piña"
$ws.Cells.Item(399, 11).Value = "Caramelo"
$ws.Cells.Item(399, 12).Value = "Segunda"
$ws.Cells.Item(399, 13).Value = 100
$ws.Cells.Item(399, 14).Value = 25000
$ws.Cells.Item(399, 15).Value = 25000
$ws.Cells.Item(399, 16).Value = 25000
$ws.Cells.Item(399, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(399, 18).Value = "Ecuador"
$ws.Cells.Item(399, 19).Value = 1786
$ws.Cells.Item(399, 20).Value = 14
